# 3dunet session annotations.xlsx - add new session log row "240204-3" (checkpoint chpt-240204-3)
# as row 170, shifting the existing rows 170:178 down to 171:179. Also clean up the
# "error boolean" column (BD) for rows 166-169 from the placeholder text "TBD" to the
# normal numeric 0 used throughout the rest of the table, and re-point the frozen-pane /
# selection to the new bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a brand-new row at 170 (pushes 170..178 down to 171..179) ------------
$ws.Rows(170).Insert()

# The inserted row inherits some per-column formatting automatically; pull the exact
# styling used by the analogous checkpoint row (169, "chpt-240204-2") for the columns
# where the auto-inherited style doesn't already match it.
$ws.Range("G169").Copy()
$ws.Range("G170").PasteSpecial(-4122)
$ws.Range("AH169").Copy()
$ws.Range("AH170").PasteSpecial(-4122)
$ws.Range("AP169").Copy()
$ws.Range("AP170").PasteSpecial(-4122)
$ws.Range("AT169").Copy()
$ws.Range("AT170").PasteSpecial(-4122)
$ws.Range("AX169").Copy()
$ws.Range("AX170").PasteSpecial(-4122)
$ws.Range("J169").Copy()
$ws.Range("J170").PasteSpecial(-4122)
$ws.Range("V169").Copy()
$ws.Range("V170").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Cells that the insert leaves behind as residual styled-but-empty cells, and that the
# new row does not actually use -- drop them entirely.
$ws.Range("U170").Clear()
$ws.Range("AM170").Clear()
$ws.Range("AY170").Clear()
$ws.Range("AZ170").Clear()
$ws.Range("BA170").Clear()

# --- 2) Fill in the new session-log row (chpt-240204-3) -----------------------------
$ws.Range("A170").Value = "240204-3"
$ws.Range("B170").Value = "autofluo eye, nuclei model type, val=id02(?), test=07"
$ws.Range("C170").Value = "dataset10.b.3"
$ws.Range("D170").Value = "train3dunet 1.8.2"
$ws.Range("E170").Value = "ResidualUNet3D"
$ws.Range("F170").Value = "3DUnet_lightsheet_nuclei"
$ws.Range("G170").Value = "Vary validation and test sample compared to chpt-240204-0, chpt-240204-1 and chpt-240204-2"
$ws.Range("H170").Value = "expect model of equal predictive power"
$ws.Range("I170").Value = "TBD"
$ws.Range("J170").Value = 1
$ws.Range("V170").Value = 1

$ws.Range("W170").Value = 6
$ws.Range("X170").Value = 5
$ws.Range("Y170").Value = 1
$ws.Range("Z170").Value = 1
$ws.Range("AA170").Formula = "=X170+Y170"
$ws.Range("AB170").Value = 6
$ws.Range("AC170").Value = 3

$ws.Range("AE170").Value = "uint16"

$ws.Range("AH170").Value = "uint8"
$ws.Range("AI170").Value = 78075
$ws.Range("AJ170").Value = 2977
$ws.Range("AK170").Formula = "=AI170+AJ170"
$ws.Range("AL170").Formula = "= 1508.06553301511 + 0.00210606006752809 * (AQ170*AR170*AS170) * (AA170 / 5) + 441"

$ws.Range("AN170").Value = 173
$ws.Range("AO170").Value = 743
$ws.Range("AP170").Value = 435
$ws.Range("AQ170").Value = 133
$ws.Range("AR170").Value = 720
$ws.Range("AS170").Value = 300
$ws.Range("AT170").Value = "yes"
$ws.Range("AU170").Formula = "= _xlfn.FLOOR.MATH((AN170 - AQ170) / 2)"
$ws.Range("AV170").Formula = "= _xlfn.FLOOR.MATH((AO170 - AR170) / 2)"
$ws.Range("AW170").Formula = "= _xlfn.FLOOR.MATH((AP170 - AS170) / 2)"
$ws.Range("AX170").Value = "yes"

$ws.Range("BB170").Value = "patch = same as for model comparison chpt-240203-5"
$ws.Range("BC170").Value = "stride = same as for model comparison chpt-240203-5"
$ws.Range("BD170").Value = 0

# --- 3) Rows 166-169: the "error boolean" column had been left as the "TBD" text ----
#        placeholder; flip it to the normal numeric 0 (no error) used elsewhere.
$ws.Range("BD166").Value = 0
$ws.Range("BD167").Value = 0
$ws.Range("BD168").Value = 0
$ws.Range("BD169").Value = 0

# --- 4) Update the view: frozen pane top-left cell & the active selection ----------
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("J154").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("A171").Select()
